$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A/B column values for rows 1..18 (order matches target sheet)
$data = @(
    @("dog", "one"),
    @("dog", "one"),
    @("dog", "one"),
    @("dog", "one"),
    @("dog", "one"),
    @("Cat", "two"),
    @("cat", "two"),
    @("cat", "two"),
    @("dog", "two"),
    @("dog", "two"),
    @("dog", "two"),
    @("dog", "two "),
    @("dog", "two "),
    @("dog", "two "),
    @("dog", "two"),
    @("dog", "two "),
    @("dog", "two "),
    @("john", "two")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $pair = $data[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

$ws.Range("B18").Select()
